$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8' date='2020-03-11 23:55:18'"
$ws.Range("A2").Value = "!!ObjTables type='Data' tableFormat='row' id='Transaction' name='Transaction' description='Stores transactions' date='2020-03-11 23:55:18' objTablesVersion='0.0.8'"

$ws.Protect()
